# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This updates the First_Noticeable_Increase_Index (C),
# First_Noticeable_Increase_Cumulative_Value (E) and the derived
# Pulse_Width (G = Point_Exceeds_Index - First_Noticeable_Increase_Index)
# columns on the four "Step3_DataPts_*" sheets.

$wb = $excel.ActiveWorkbook

# New First_Noticeable_Increase_Index (column C) and
# First_Noticeable_Increase_Cumulative_Value (column E) values, shared by
# every Step3_DataPts_* sheet (rows 2-6).
$newC = @{ 2 = 87; 3 = 87; 4 = 88; 5 = 88; 6 = 87 }
$newE = @{
    2 = 0.0006972086426067162
    3 = 0.01717146190835171
    4 = 0.01819074046114215
    5 = 0.02555204893951814
    6 = 0.02508668387264054
}

$sheetNames = @("Step3_DataPts_0.5", "Step3_DataPts_0.7", "Step3_DataPts_0.8", "Step3_DataPts_0.9")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in 2..6) {
        $ws.Cells.Item($row, 3).Value = $newC[$row]
        $ws.Cells.Item($row, 5).Value = $newE[$row]

        $pointExceedsIndex = $ws.Cells.Item($row, 4).Value2
        $ws.Cells.Item($row, 7).Value = $pointExceedsIndex - $newC[$row]
    }
}
